$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '50.992.76'
$ws.Range("E2").Value = '  +0.04%  '

# Row 3
$ws.Range("D3").Value = '2.952.09'
$ws.Range("E3").Value = '  +0.89%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '380.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.34%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.17%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.541'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.04%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.64%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.52%  '

# Row 11
$ws.Range("E11").Value = '  -0.48%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0849'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.62%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.70%  '

# Row 14
$ws.Range("D14").Value = '3.420.82'
$ws.Range("E14").Value = '  +0.80%  '

# Row 15
$ws.Range("B15").Value = 'Uniswap'
$ws.Range("C15").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '12.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +74.45%  '

# Row 16
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.56%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.955.07'
$ws.Range("E17").Value = '  +1.57%  '

# Row 18
$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.10%  '

# Row 19
$ws.Range("D19").Value = '51.007.72'
$ws.Range("E19").Value = '  +0.26%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.57%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0953'
$ws.Range("E22").Value = '  +0.42%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +18.36%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.17%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '266.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.90%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.27%  '

# Row 27
$ws.Range("E27").Value = '  +0.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.82%  '

# Row 29
$ws.Range("E29").Value = '  -0.78%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.36%  '

# Row 31
$ws.Range("E31").Value = '  -6.71%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.50%  '

# Row 33
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.22%  '

# Row 34
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.99%  '

# Row 35
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.93'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0433'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.63%  '

# Row 37
$ws.Range("E37").Value = '  +0.12%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.46%  '

# Row 39
$ws.Range("E39").Value = '  +1.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.65%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.82'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.88%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.67%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.30%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.69%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.12%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.76%  '

# Row 47
$ws.Range("D47").Value = '2.023.99'
$ws.Range("E47").Value = '  +1.22%  '

# Row 48
$ws.Range("E48").Value = '  -2.52%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.255'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.26%  '

# Row 50
$ws.Range("E50").Value = '  -7.66%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.34'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.53%  '
